{"js": "// Welcome.docx (German) edit:\n// Remove the trailing clause \", until Smartcash reaches a considerable\n// market cap\" from the ASIC/mining paragraph, leaving the sentence ending\n// at \"...created for quite some time.\"\n\nconst body = context.document.body;\n\nconst oldClause = \", until Smartcash reaches a considerable market cap.\";\nconst newEnding = \".\";\n\nconst results = body.search(oldClause, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target text not found: \"' + oldClause + '\"');\n}\n\n// Replace the matched range's text (keeps the paragraph's single run and\n// its run-level formatting, it only trims the trailing clause).\nresults.items[0].insertText(newEnding, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Welcome.docx (German) edit:\n# Remove the trailing clause \", until Smartcash reaches a considerable\n# market cap\" from the ASIC/mining paragraph, leaving the sentence ending\n# at \"...created for quite some time.\"\n\n$d = $word.ActiveDocument\n\n$oldClause = \", until Smartcash reaches a considerable market cap.\"\n$newEnding = \".\"\n\n$matchCount = 0\nforeach ($p in $d.Paragraphs) {\n    $rng = $p.Range\n    $text = $rng.Text\n    if ($text -like \"*$oldClause*\") {\n        $matchCount++\n\n        # Paragraph.Range.Text includes the trailing paragraph-mark (CR,\n        # char 13); strip it before editing so re-assigning .Text doesn't\n        # insert an extra paragraph break.\n        if ($text.Length -gt 0 -and [int][char]$text.Substring($text.Length - 1, 1) -eq 13) {\n            $text = $text.Substring(0, $text.Length - 1)\n        }\n\n        $rng.Text = $text.Replace($oldClause, $newEnding)\n    }\n}\n\nif ($matchCount -eq 0) {\n    throw \"Target text not found: $oldClause\"\n}\n"}
